# Auto-generated Excel COM-interop edit script
# Applies updated crypto price/volume data to Sheet1 per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.877.35"
$ws.Range("E2").Value = "'  -0.35%  "

# Row 3
$ws.Range("D3").Value = "'2.347.84"
$ws.Range("E3").Value = "'  -0.10%  "

# Row 4
$ws.Range("E4").Value = "'  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'239.36"
$ws.Range("E5").Value = "'  +0.14%  "

# Row 6
$ws.Range("D6").Value = "'0.666"
$ws.Range("E6").Value = "'  -1.68%  "

# Row 7
$ws.Range("D7").Value = "'73.20"
$ws.Range("E7").Value = "'  -1.02%  "

# Row 8
$ws.Range("E8").Value = "'  +0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.597"
$ws.Range("E9").Value = "'  +0.50%  "

# Row 10
$ws.Range("B10").Value = "'OKB"
$ws.Range("C10").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'61.92"
$ws.Range("E10").Value = "'  +8.33%  "

# Row 11
$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = "'  +0.40%  "

# Row 12
$ws.Range("D12").Value = "'33.54"
$ws.Range("E12").Value = "'  +4.92%  "

# Row 13
$ws.Range("E13").Value = "'  +0.78%  "

# Row 14
$ws.Range("D14").Value = "'7.20"
$ws.Range("E14").Value = "'  +0.38%  "

# Row 15
$ws.Range("D15").Value = "'16.11"
$ws.Range("E15").Value = "'  -2.27%  "

# Row 16
$ws.Range("D16").Value = "'0.902"
$ws.Range("E16").Value = "'  +0.54%  "

# Row 17
$ws.Range("D17").Value = "'2.346.15"
$ws.Range("E17").Value = "'  -0.35%  "

# Row 18
$ws.Range("D18").Value = "'43.783.23"
$ws.Range("E18").Value = "'  -0.30%  "

# Row 19
$ws.Range("E19").Value = "'  -0.05%  "

# Row 20
$ws.Range("B20").Value = "'Litecoin"
$ws.Range("C20").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'77.80"
$ws.Range("E20").Value = "'  +1.35%  "

# Row 21
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.57"
$ws.Range("E21").Value = "'  -1.12%  "

# Row 22
$ws.Range("D22").Value = "'252.48"
$ws.Range("E22").Value = "'  -1.39%  "

# Row 23
$ws.Range("E23").Value = "'  +3.02%  "

# Row 24
$ws.Range("E24").Value = "'  +0.05%  "

# Row 25
$ws.Range("E25").Value = "'  -3.65%  "

# Row 26
$ws.Range("E26").Value = "'  -0.34%  "

# Row 27
$ws.Range("D27").Value = "'10.39"
$ws.Range("E27").Value = "'  -2.43%  "

# Row 28
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "'  +0.73%  "

# Row 29
$ws.Range("D29").Value = "'175.62"
$ws.Range("E29").Value = "'  +0.58%  "

# Row 30
$ws.Range("D30").Value = "'22.17"
$ws.Range("E30").Value = "'  -2.40%  "

# Row 31
$ws.Range("E31").Value = "'  +0.78%  "

# Row 32
$ws.Range("E32").Value = "'  -2.21%  "

# Row 33
$ws.Range("D33").Value = "'0.0739"
$ws.Range("E33").Value = "'  -2.38%  "

# Row 34
$ws.Range("D34").Value = "'5.04"
$ws.Range("E34").Value = "'  -4.34%  "

# Row 35
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "'  -0.39%  "

# Row 36
$ws.Range("D36").Value = "'3.74"
$ws.Range("E36").Value = "'  +0.86%  "

# Row 37
$ws.Range("B37").Value = "'THORChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.41"
$ws.Range("E37").Value = "'  +1.48%  "

# Row 38
$ws.Range("B38").Value = "'LidoDAOToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.39"
$ws.Range("E38").Value = "'  +1.46%  "

# Row 39
$ws.Range("D39").Value = "'0.0271"
$ws.Range("E39").Value = "'  -3.52%  "

# Row 40
$ws.Range("D40").Value = "'5.40"
$ws.Range("E40").Value = "'  +15.96%  "

# Row 41
$ws.Range("D41").Value = "'64.82"
$ws.Range("E41").Value = "'  +12.44%  "

# Row 42
$ws.Range("D42").Value = "'19.45"
$ws.Range("E42").Value = "'  +2.18%  "

# Row 43
$ws.Range("D43").Value = "'9.09"
$ws.Range("E43").Value = "'  +0.72%  "

# Row 44
$ws.Range("D44").Value = "'0.106"
$ws.Range("E44").Value = "'  -2.39%  "

# Row 45
$ws.Range("D45").Value = "'0.199"
$ws.Range("E45").Value = "'  -1.07%  "

# Row 46
$ws.Range("E46").Value = "'  +0.14%  "

# Row 47
$ws.Range("E47").Value = "'  -1.92%  "

# Row 48
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "'  -2.27%  "

# Row 49
$ws.Range("E49").Value = "'  -1.84%  "

# Row 50
$ws.Range("B50").Value = "'TerraClassic"
$ws.Range("C50").Value = "'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").Value = "'0.000214"
$ws.Range("E50").Value = "'  +17.47%  "

# Row 51
$ws.Range("B51").Value = "'Aave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'97.54"
$ws.Range("E51").Value = "'  -2.31%  "
